# Updated cryptos list on Thu Jun 27 21:07:49 UTC 2024 with GitHub Actions
#
# Refreshes the Price (col D) and Volume(1h) (col E) figures scraped for
# each coin row. Two coin pairs also swapped rank this run, so their
# Coin name (B) and Link (C) cells are rewritten along with D/E:
#   - rows 26/27: LEO <-> Dai
#   - rows 44/45: FirstDigitalUSD <-> Filecoin
#
# Price cells are leading-apostrophe prefixed so Excel keeps them as text
# (matching the source data, which uses literal "."-grouped strings like
# "61.503.50" and fixed-precision strings like "1.00" rather than numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.503.50"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "'3.445.02"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'579.73"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'149.42"
$ws.Range("E6").Value = "  +8.84%  "
$ws.Range("D7").Value = "'3.445.41"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'7.85"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'4.033.32"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "'28.06"
$ws.Range("E14").Value = "  +6.50%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "'3.447.05"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "'61.577.01"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +8.38%  "
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "'9.50"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'389.05"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "'3.586.47"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "'72.87"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.77"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "'0.182"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").Value = "'7.79"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'1.52"
$ws.Range("E32").Value = "  -14.77%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'24.05"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'5.32"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").Value = "'7.09"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "'166.59"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'0.0795"
$ws.Range("E41").Value = "  +4.39%  "
$ws.Range("D42").Value = "'26.48"
$ws.Range("E42").Value = "  +9.31%  "
$ws.Range("D43").Value = "'0.795"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.51"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'42.33"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'2.617.36"
$ws.Range("E48").Value = "  +5.28%  "
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("D50").Value = "'7.08"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("D51").Value = "'23.21"
$ws.Range("E51").Value = "  -1.18%  "
